$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 12 de Octubre de 2020 a las 01:27"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 7990313
$ws.Cells.Item(4, 3).Value = 40993
$ws.Cells.Item(4, 4).Value = 5126070
$ws.Cells.Item(4, 5).Value = 2644557
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 316
$ws.Cells.Item(4, 8).Value = 219686

# Row 8: Colombia
$ws.Cells.Item(8, 1).Value = "Colombia"
$ws.Cells.Item(8, 2).Value = 911316
$ws.Cells.Item(8, 3).Value = 8569
$ws.Cells.Item(8, 4).Value = 789787
$ws.Cells.Item(8, 5).Value = 93695
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 174
$ws.Cells.Item(8, 8).Value = 27834

# Row 9: Argentina
$ws.Cells.Item(9, 1).Value = "Argentina"
$ws.Cells.Item(9, 2).Value = 894206
$ws.Cells.Item(9, 3).Value = 10324
$ws.Cells.Item(9, 4).Value = 721380
$ws.Cells.Item(9, 5).Value = 148958
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 287
$ws.Cells.Item(9, 8).Value = 23868

# Row 10: España
$ws.Cells.Item(10, 1).Value = "España"
$ws.Cells.Item(10, 2).Value = 890367
$ws.Cells.Item(10, 3).Value = 0
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 32929

# Row 36: Catar
$ws.Cells.Item(36, 1).Value = "Catar"
$ws.Cells.Item(36, 2).Value = 127985
$ws.Cells.Item(36, 3).Value = 207
$ws.Cells.Item(36, 4).Value = 124978
$ws.Cells.Item(36, 5).Value = 2787
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 1
$ws.Cells.Item(36, 8).Value = 220

# Row 38: Panama
$ws.Cells.Item(38, 1).Value = "Panama"
$ws.Cells.Item(38, 2).Value = 120313
$ws.Cells.Item(38, 3).Value = 647
$ws.Cells.Item(38, 4).Value = 96164
$ws.Cells.Item(38, 5).Value = 21658
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 9
$ws.Cells.Item(38, 8).Value = 2491

# Row 40: Chequia
$ws.Cells.Item(40, 1).Value = "Chequia"
$ws.Cells.Item(40, 2).Value = 117110
$ws.Cells.Item(40, 3).Value = 3105
$ws.Cells.Item(40, 4).Value = 54980
$ws.Cells.Item(40, 5).Value = 61143
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 39
$ws.Cells.Item(40, 8).Value = 987

# Row 46: Egipto
$ws.Cells.Item(46, 1).Value = "Egipto"
$ws.Cells.Item(46, 2).Value = 104516
$ws.Cells.Item(46, 3).Value = 129
$ws.Cells.Item(46, 4).Value = 97688
$ws.Cells.Item(46, 5).Value = 776
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(46, 7).Value = 12
$ws.Cells.Item(46, 8).Value = 6052

# Row 57: Barein
$ws.Cells.Item(57, 1).Value = "Barein"
$ws.Cells.Item(57, 2).Value = 75614
$ws.Cells.Item(57, 3).Value = 327
$ws.Cells.Item(57, 4).Value = 71249
$ws.Cells.Item(57, 5).Value = 4090
$ws.Cells.Item(57, 6).Value = 0
$ws.Cells.Item(57, 7).Value = 2
$ws.Cells.Item(57, 8).Value = 275

# Row 61: Nigeria
$ws.Cells.Item(61, 1).Value = "Nigeria"
$ws.Cells.Item(61, 2).Value = 60266
$ws.Cells.Item(61, 3).Value = 163
$ws.Cells.Item(61, 4).Value = 51735
$ws.Cells.Item(61, 5).Value = 7416
$ws.Cells.Item(61, 6).Value = 0
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 1115

# Row 67: Paraguay
$ws.Cells.Item(67, 1).Value = "Paraguay"
$ws.Cells.Item(67, 2).Value = 49675
$ws.Cells.Item(67, 3).Value = 697
$ws.Cells.Item(67, 4).Value = 32090
$ws.Cells.Item(67, 5).Value = 16508
$ws.Cells.Item(67, 6).Value = 0
$ws.Cells.Item(67, 7).Value = 12
$ws.Cells.Item(67, 8).Value = 1077

# Row 68: Kirguistan
$ws.Cells.Item(68, 1).Value = "Kirguistan"
$ws.Cells.Item(68, 2).Value = 49230
$ws.Cells.Item(68, 3).Value = 306
$ws.Cells.Item(68, 4).Value = 44227
$ws.Cells.Item(68, 5).Value = 3918
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 3
$ws.Cells.Item(68, 8).Value = 1085

# Row 78: Tunez
$ws.Cells.Item(78, 1).Value = "Tunez"
$ws.Cells.Item(78, 2).Value = 32556
$ws.Cells.Item(78, 3).Value = 1297
$ws.Cells.Item(78, 4).Value = 5032
$ws.Cells.Item(78, 5).Value = 27046
$ws.Cells.Item(78, 6).Value = 0
$ws.Cells.Item(78, 7).Value = 22
$ws.Cells.Item(78, 8).Value = 478

# Row 79: Dinamarca
$ws.Cells.Item(79, 1).Value = "Dinamarca"
$ws.Cells.Item(79, 2).Value = 32422
$ws.Cells.Item(79, 3).Value = 340
$ws.Cells.Item(79, 4).Value = 26380
$ws.Cells.Item(79, 5).Value = 5373
$ws.Cells.Item(79, 6).Value = 0
$ws.Cells.Item(79, 7).Value = 2
$ws.Cells.Item(79, 8).Value = 669

# Row 95: Noruega
$ws.Cells.Item(95, 1).Value = "Noruega"
$ws.Cells.Item(95, 2).Value = 15524
$ws.Cells.Item(95, 3).Value = 58
$ws.Cells.Item(95, 4).Value = 11863
$ws.Cells.Item(95, 5).Value = 3386
$ws.Cells.Item(95, 6).Value = 0
$ws.Cells.Item(95, 7).Value = 0
$ws.Cells.Item(95, 8).Value = 275

# Row 112: Haiti
$ws.Cells.Item(112, 1).Value = "Haiti"
$ws.Cells.Item(112, 2).Value = 8882
$ws.Cells.Item(112, 3).Value = 22
$ws.Cells.Item(112, 4).Value = 7104
$ws.Cells.Item(112, 5).Value = 1548
$ws.Cells.Item(112, 6).Value = 0
$ws.Cells.Item(112, 7).Value = 0
$ws.Cells.Item(112, 8).Value = 230

# Row 124: Suazilandia
$ws.Cells.Item(124, 1).Value = "Suazilandia"
$ws.Cells.Item(124, 2).Value = 5669
$ws.Cells.Item(124, 3).Value = 9
$ws.Cells.Item(124, 4).Value = 5310
$ws.Cells.Item(124, 5).Value = 246
$ws.Cells.Item(124, 6).Value = 0
$ws.Cells.Item(124, 7).Value = 0
$ws.Cells.Item(124, 8).Value = 113

# Row 157: Uruguay
$ws.Cells.Item(157, 1).Value = "Uruguay"
$ws.Cells.Item(157, 2).Value = 2294
$ws.Cells.Item(157, 3).Value = 26
$ws.Cells.Item(157, 4).Value = 1942
$ws.Cells.Item(157, 5).Value = 302
$ws.Cells.Item(157, 6).Value = 0
$ws.Cells.Item(157, 7).Value = 0
$ws.Cells.Item(157, 8).Value = 50

# Row 158: Burkina Faso
$ws.Cells.Item(158, 1).Value = "Burkina Faso"
$ws.Cells.Item(158, 2).Value = 2271
$ws.Cells.Item(158, 3).Value = 0
$ws.Cells.Item(158, 4).Value = 1542
$ws.Cells.Item(158, 5).Value = 668
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 0
$ws.Cells.Item(158, 8).Value = 61
